# Natmi following Dr Hou advice:
# Re-derive the Wnt5a-Lrp5 LR-pair sheet with TWO sending clusters
# (ECs and FAPs) instead of just FAPs. Rows 2-6 (previously FAPs -> {ECs,
# FAPs, M1, M2, sCs}) are rewritten as ECs -> {ECs, FAPs, M1, M2, sCs}
# with recomputed stats, and rows 7-11 are appended for
# FAPs -> {ECs, FAPs, M1, M2, sCs}.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Wnt5a", "Lrp5", "ECs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 17.82379233333333, 53.471377, 0.4233776263711466, 0.4233776263711467, 2.850030335364111, 25.650273018277, 0.01121007400597485, 0.01121007400597486),
  @("ECs", "Wnt5a", "Lrp5", "FAPs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 10.453073, 31.359219, 0.2482971722436279, 0.2482971722436279, 1.671449857057667, 15.043048713519, 0.006574342863090524, 0.006574342863090525),
  @("ECs", "Wnt5a", "Lrp5", "M1", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 4.27602, 12.82806, 0.1015704830969034, 0.1015704830969034, 0.68373702334, 6.15363321006, 0.002689354754284443, 0.002689354754284444),
  @("ECs", "Wnt5a", "Lrp5", "M2", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 6.530620666666667, 19.591862, 0.1551251621919343, 0.1551251621919343, 1.044248421473556, 9.398235793262, 0.004107360521776848, 0.004107360521776848),
  @("ECs", "Wnt5a", "Lrp5", "sCs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 3.015535666666667, 9.046607, 0.07162955609638778, 0.07162955609638778, 0.4821851582785556, 4.339666424507, 0.001896587289550635, 0.001896587289550635),
  @("FAPs", "Wnt5a", "Lrp5", "ECs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 17.82379233333333, 53.471377, 0.4233776263711466, 0.4233776263711467, 104.7887843441013, 943.0990590969121, 0.4121675523651717, 0.4121675523651718),
  @("FAPs", "Wnt5a", "Lrp5", "FAPs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 10.453073, 31.359219, 0.2482971722436279, 0.2482971722436279, 61.455205034096, 553.096845306864, 0.2417228293805374, 0.2417228293805374),
  @("FAPs", "Wnt5a", "Lrp5", "M1", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 4.27602, 12.82806, 0.1015704830969034, 0.1015704830969034, 25.13937153504, 226.25434381536, 0.09888112834261899, 0.09888112834261902),
  @("FAPs", "Wnt5a", "Lrp5", "M2", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 6.530620666666667, 19.591862, 0.1551251621919343, 0.1551251621919343, 38.39451155367467, 345.550603983072, 0.1510178016701575, 0.1510178016701574),
  @("FAPs", "Wnt5a", "Lrp5", "sCs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 3.015535666666667, 9.046607, 0.07162955609638778, 0.07162955609638778, 17.72879254575467, 159.559132911792, 0.06973296880683715, 0.06973296880683716),
)

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
# D Target cluster, E..T the numeric LR statistics.
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  $r = $i + 2
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}
